$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Fecha) - serial date values shuffled across rows 2-9
$ws.Range("D2").Value = 44323
$ws.Range("D3").Value = 44309
$ws.Range("D4").Value = 44322
$ws.Range("D5").Value = 44306
$ws.Range("D6").Value = 44327
$ws.Range("D7").Value = 44313
$ws.Range("D8").Value = 44302
$ws.Range("D9").Value = 44330

# Column M (Volumen) - values shuffled across rows 2,4,6,7,8
$ws.Range("M2").Value = 80
$ws.Range("M4").Value = 60
$ws.Range("M6").Value = 60
$ws.Range("M7").Value = 120
$ws.Range("M8").Value = 80

# Rows 3 and 6 - Unidad de comercialización (Q), Precio $/Kg (S) and Kg/unidad (T) swapped
$ws.Range("Q3").Value = "$/caja 14 kilos granel"
$ws.Range("S3").Value = 821
$ws.Range("T3").Value = 14

$ws.Range("Q6").Value = "$/caja 10 kilos empedrada"
$ws.Range("S6").Value = 11500
$ws.Range("T6").Value = 1
